$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated loading_percent results for the 380 kV case (Case_3_113).
# Columns B:D, F:I and M:O are rewritten per data row (2-25); E, J, K, L are untouched (stay 0).

$ws.Cells.Item(2, 2).Value = 9.741583706518572
$ws.Cells.Item(2, 3).Value = 8.371277728407208
$ws.Cells.Item(2, 4).Value = 3.510785840105738
$ws.Cells.Item(2, 6).Value = 17.55709551358513
$ws.Cells.Item(2, 7).Value = 17.48625437641607
$ws.Cells.Item(2, 8).Value = 11.22903656869528
$ws.Cells.Item(2, 9).Value = 15.68704182894419
$ws.Cells.Item(2, 13).Value = 19.96626076887188
$ws.Cells.Item(2, 14).Value = 15.92661311942672
$ws.Cells.Item(2, 15).Value = 15.58854796573555
$ws.Cells.Item(3, 2).Value = 9.225687150525037
$ws.Cells.Item(3, 3).Value = 8.142480099828742
$ws.Cells.Item(3, 4).Value = 3.417341587081574
$ws.Cells.Item(3, 6).Value = 17.50238079788274
$ws.Cells.Item(3, 7).Value = 17.34255310698098
$ws.Cells.Item(3, 8).Value = 11.25890365256534
$ws.Cells.Item(3, 9).Value = 15.76500195294433
$ws.Cells.Item(3, 13).Value = 19.35638236821036
$ws.Cells.Item(3, 14).Value = 15.97917663091778
$ws.Cells.Item(3, 15).Value = 15.60998301096311
$ws.Cells.Item(4, 2).Value = 8.891830553975391
$ws.Cells.Item(4, 3).Value = 7.997621720431522
$ws.Cells.Item(4, 4).Value = 3.357954955879177
$ws.Cells.Item(4, 6).Value = 17.47438116032757
$ws.Cells.Item(4, 7).Value = 17.26155192603204
$ws.Cells.Item(4, 8).Value = 11.279265940445
$ws.Cells.Item(4, 9).Value = 15.81646491323354
$ws.Cells.Item(4, 13).Value = 18.9800296979327
$ws.Cells.Item(4, 14).Value = 16.01316683363613
$ws.Cells.Item(4, 15).Value = 15.62730892758882
$ws.Cells.Item(5, 2).Value = 8.751559694818999
$ws.Cells.Item(5, 3).Value = 7.937551454779205
$ws.Cells.Item(5, 4).Value = 3.333268822993516
$ws.Cells.Item(5, 6).Value = 17.46438749426902
$ws.Cells.Item(5, 7).Value = 17.23040048897764
$ws.Cells.Item(5, 8).Value = 11.28807212883149
$ws.Cells.Item(5, 9).Value = 15.83833880926653
$ws.Cells.Item(5, 13).Value = 18.82644902742003
$ws.Cells.Item(5, 14).Value = 16.02745102908344
$ws.Cells.Item(5, 15).Value = 15.63541435319821
$ws.Cells.Item(6, 2).Value = 8.728015419522132
$ws.Cells.Item(6, 3).Value = 7.927515937304348
$ws.Cells.Item(6, 4).Value = 3.329141002521772
$ws.Cells.Item(6, 6).Value = 17.46281383823455
$ws.Cells.Item(6, 7).Value = 17.22534102718456
$ws.Cells.Item(6, 8).Value = 11.28956507976662
$ws.Cells.Item(6, 9).Value = 15.84202540190952
$ws.Cells.Item(6, 13).Value = 18.8009420166901
$ws.Cells.Item(6, 14).Value = 16.02984909723813
$ws.Cells.Item(6, 15).Value = 15.63682329135251
$ws.Cells.Item(7, 2).Value = 8.889955791817137
$ws.Cells.Item(7, 3).Value = 7.996815715874078
$ws.Cells.Item(7, 4).Value = 3.357623968364376
$ws.Cells.Item(7, 6).Value = 17.47424063629123
$ws.Cells.Item(7, 7).Value = 17.26112423945605
$ws.Cells.Item(7, 8).Value = 11.27938264603512
$ws.Cells.Item(7, 9).Value = 15.81675626103888
$ws.Cells.Item(7, 13).Value = 18.9779589679748
$ws.Cells.Item(7, 14).Value = 16.01335772062422
$ws.Cells.Item(7, 15).Value = 15.62741401250559
$ws.Cells.Item(8, 2).Value = 9.567291728810677
$ws.Cells.Item(8, 3).Value = 8.293331621183773
$ws.Cells.Item(8, 4).Value = 3.478995367039199
$ws.Cells.Item(8, 6).Value = 17.53707477471432
$ws.Cells.Item(8, 7).Value = 17.43522976238365
$ws.Cells.Item(8, 8).Value = 11.23891441185313
$ws.Cells.Item(8, 9).Value = 15.71317528156395
$ws.Cells.Item(8, 13).Value = 19.75651757742983
$ws.Cells.Item(8, 14).Value = 15.94438171282279
$ws.Cells.Item(8, 15).Value = 15.59507314767917
$ws.Cells.Item(9, 2).Value = 10.75727954343059
$ws.Cells.Item(9, 3).Value = 8.837639506270888
$ws.Cells.Item(9, 4).Value = 3.700238607488843
$ws.Cells.Item(9, 6).Value = 17.70418603239189
$ws.Cells.Item(9, 7).Value = 17.83211526639604
$ws.Cells.Item(9, 8).Value = 11.17563558233223
$ws.Cells.Item(9, 9).Value = 15.53865576041831
$ws.Cells.Item(9, 13).Value = 21.25764281709821
$ws.Cells.Item(9, 14).Value = 15.82267151432452
$ws.Cells.Item(9, 15).Value = 15.56478068287443
$ws.Cells.Item(10, 2).Value = 11.54459911033798
$ws.Cells.Item(10, 3).Value = 9.211952857337781
$ws.Cells.Item(10, 4).Value = 3.851598391713686
$ws.Cells.Item(10, 6).Value = 17.8529195316437
$ws.Cells.Item(10, 7).Value = 18.15467291242426
$ws.Cells.Item(10, 8).Value = 11.13897469957832
$ws.Cells.Item(10, 9).Value = 15.42797914398383
$ws.Cells.Item(10, 13).Value = 22.33179238043865
$ws.Cells.Item(10, 14).Value = 15.74142306901773
$ws.Cells.Item(10, 15).Value = 15.56280514888808
$ws.Cells.Item(11, 2).Value = 11.88351659220155
$ws.Cells.Item(11, 3).Value = 9.376150451824564
$ws.Cells.Item(11, 4).Value = 3.917852257888853
$ws.Cells.Item(11, 6).Value = 17.92600988949764
$ws.Cells.Item(11, 7).Value = 18.30746327351255
$ws.Cells.Item(11, 8).Value = 11.12443633376326
$ws.Cells.Item(11, 9).Value = 15.38146231613817
$ws.Cells.Item(11, 13).Value = 22.81179960363785
$ws.Cells.Item(11, 14).Value = 15.70621660638424
$ws.Cells.Item(11, 15).Value = 15.56631775140987
$ws.Cells.Item(12, 2).Value = 12.0090676101217
$ws.Cells.Item(12, 3).Value = 9.437414948132004
$ws.Cells.Item(12, 4).Value = 3.942554210596748
$ws.Cells.Item(12, 6).Value = 17.95444810967084
$ws.Cells.Item(12, 7).Value = 18.36613217694134
$ws.Cells.Item(12, 8).Value = 11.11923894182841
$ws.Cells.Item(12, 9).Value = 15.364400517183
$ws.Cells.Item(12, 13).Value = 22.99214861276137
$ws.Cells.Item(12, 14).Value = 15.69313560726183
$ws.Cells.Item(12, 15).Value = 15.56828191141002
$ws.Cells.Item(13, 2).Value = 11.98215233155375
$ws.Cells.Item(13, 3).Value = 9.424261758174572
$ws.Cells.Item(13, 4).Value = 3.93725160581608
$ws.Cells.Item(13, 6).Value = 17.94828994913664
$ws.Cells.Item(13, 7).Value = 18.35346171253959
$ws.Cells.Item(13, 8).Value = 11.12034459013278
$ws.Cells.Item(13, 9).Value = 15.36805044408755
$ws.Cells.Item(13, 13).Value = 22.95337300185998
$ws.Cells.Item(13, 14).Value = 15.69594169476633
$ws.Cells.Item(13, 15).Value = 15.56783070753273
$ws.Cells.Item(14, 2).Value = 11.89390177153785
$ws.Cells.Item(14, 3).Value = 9.381209239326784
$ws.Cells.Item(14, 4).Value = 3.919892327233093
$ws.Cells.Item(14, 6).Value = 17.92833440027221
$ws.Cells.Item(14, 7).Value = 18.31227410831891
$ws.Cells.Item(14, 8).Value = 11.12400256684778
$ws.Cells.Item(14, 9).Value = 15.38004753318911
$ws.Cells.Item(14, 13).Value = 22.8266664286378
$ws.Cells.Item(14, 14).Value = 15.70513540310412
$ws.Cells.Item(14, 15).Value = 15.56646664128362
$ws.Cells.Item(15, 2).Value = 11.83948186820536
$ws.Cells.Item(15, 3).Value = 9.354718255040803
$ws.Cells.Item(15, 4).Value = 3.909208512336237
$ws.Cells.Item(15, 6).Value = 17.91620945512926
$ws.Cells.Item(15, 7).Value = 18.28714923847286
$ws.Cells.Item(15, 8).Value = 11.12628330138486
$ws.Cells.Item(15, 9).Value = 15.38746819706966
$ws.Cells.Item(15, 13).Value = 22.74886516815942
$ws.Cells.Item(15, 14).Value = 15.71079945490838
$ws.Cells.Item(15, 15).Value = 15.56571365954414
$ws.Cells.Item(16, 2).Value = 11.52206080868006
$ws.Cells.Item(16, 3).Value = 9.201096368549479
$ws.Cells.Item(16, 4).Value = 3.847215075662171
$ws.Cells.Item(16, 6).Value = 17.84825047323502
$ws.Cells.Item(16, 7).Value = 18.14480445827232
$ws.Cells.Item(16, 8).Value = 11.13996789118607
$ws.Cells.Item(16, 9).Value = 15.43109640283932
$ws.Cells.Item(16, 13).Value = 22.30023347462019
$ws.Cells.Item(16, 14).Value = 15.74375907670677
$ws.Cells.Item(16, 15).Value = 15.56266434333681
$ws.Cells.Item(17, 2).Value = 11.3223874105113
$ws.Cells.Item(17, 3).Value = 9.105267941178708
$ws.Cells.Item(17, 4).Value = 3.808508390918255
$ws.Cells.Item(17, 6).Value = 17.80793673650782
$ws.Cells.Item(17, 7).Value = 18.05898948717797
$ws.Cells.Item(17, 8).Value = 11.14891106183029
$ws.Cells.Item(17, 9).Value = 15.45884359940808
$ws.Cells.Item(17, 13).Value = 22.02266966164587
$ws.Cells.Item(17, 14).Value = 15.76442702909561
$ws.Cells.Item(17, 15).Value = 15.56192346734513
$ws.Cells.Item(18, 2).Value = 11.20573073138834
$ws.Cells.Item(18, 3).Value = 9.0495798280953
$ws.Cells.Item(18, 4).Value = 3.786001252068102
$ws.Cells.Item(18, 6).Value = 17.78526183600683
$ws.Cells.Item(18, 7).Value = 18.01020487981038
$ws.Cells.Item(18, 8).Value = 11.15425623506608
$ws.Cells.Item(18, 9).Value = 15.47516346044148
$ws.Cells.Item(18, 13).Value = 21.86221962219868
$ws.Cells.Item(18, 14).Value = 15.77647984620755
$ws.Cells.Item(18, 15).Value = 15.56191260031293
$ws.Cells.Item(19, 2).Value = 11.16592275466823
$ws.Cells.Item(19, 3).Value = 9.030628133819722
$ws.Cells.Item(19, 4).Value = 3.778339219226509
$ws.Cells.Item(19, 6).Value = 17.77767312291253
$ws.Cells.Item(19, 7).Value = 17.99378763095633
$ws.Cells.Item(19, 8).Value = 11.1561005787184
$ws.Cells.Item(19, 9).Value = 15.48075091609181
$ws.Cells.Item(19, 13).Value = 21.80776187104436
$ws.Cells.Item(19, 14).Value = 15.78058912810182
$ws.Cells.Item(19, 15).Value = 15.5619802438737
$ws.Cells.Item(20, 2).Value = 11.34383048179041
$ws.Cells.Item(20, 3).Value = 9.115528345634425
$ws.Cells.Item(20, 4).Value = 3.812654143111761
$ws.Cells.Item(20, 6).Value = 17.81217530452126
$ws.Cells.Item(20, 7).Value = 18.06806567426953
$ws.Cells.Item(20, 8).Value = 11.14793820935639
$ws.Cells.Item(20, 9).Value = 15.45585254570515
$ws.Cells.Item(20, 13).Value = 22.05230122763431
$ws.Cells.Item(20, 14).Value = 15.76220980641597
$ws.Cells.Item(20, 15).Value = 15.56195935707781
$ws.Cells.Item(21, 2).Value = 11.91989896735342
$ws.Cells.Item(21, 3).Value = 9.393879895566878
$ws.Cells.Item(21, 4).Value = 3.925001766921676
$ws.Cells.Item(21, 6).Value = 17.93417535986972
$ws.Cells.Item(21, 7).Value = 18.32435039973404
$ws.Cells.Item(21, 8).Value = 11.12291976917129
$ws.Cells.Item(21, 9).Value = 15.3765086675072
$ws.Cells.Item(21, 13).Value = 22.86392307352413
$ws.Cells.Item(21, 14).Value = 15.70242818822288
$ws.Cells.Item(21, 15).Value = 15.56685009898502
$ws.Cells.Item(22, 2).Value = 12.2801258223383
$ws.Cells.Item(22, 3).Value = 9.57045825551571
$ws.Cells.Item(22, 4).Value = 3.99616688345522
$ws.Cells.Item(22, 6).Value = 18.01833160845554
$ws.Cells.Item(22, 7).Value = 18.49654241010762
$ws.Cells.Item(22, 8).Value = 11.10836401772952
$ws.Cells.Item(22, 9).Value = 15.3278779180913
$ws.Cells.Item(22, 13).Value = 23.38602183484379
$ws.Cells.Item(22, 14).Value = 15.66481943846121
$ws.Cells.Item(22, 15).Value = 15.57374145541544
$ws.Cells.Item(23, 2).Value = 12.08936062305434
$ws.Cells.Item(23, 3).Value = 9.476715695902472
$ws.Cells.Item(23, 4).Value = 3.958395521156839
$ws.Cells.Item(23, 6).Value = 17.9730183538994
$ws.Cells.Item(23, 7).Value = 18.40423109152324
$ws.Cells.Item(23, 8).Value = 11.11596831042873
$ws.Cells.Item(23, 9).Value = 15.3535372045782
$ws.Cells.Item(23, 13).Value = 23.10818563354194
$ws.Cells.Item(23, 14).Value = 15.68475858276911
$ws.Cells.Item(23, 15).Value = 15.56972556884581
$ws.Cells.Item(24, 2).Value = 11.33414186647783
$ws.Cells.Item(24, 3).Value = 9.11089146727298
$ws.Cells.Item(24, 4).Value = 3.810780639165326
$ws.Cells.Item(24, 6).Value = 17.81025748276117
$ws.Cells.Item(24, 7).Value = 18.06396060752148
$ws.Cells.Item(24, 8).Value = 11.14837740174381
$ws.Cells.Item(24, 9).Value = 15.45720365617775
$ws.Cells.Item(24, 13).Value = 22.03890750804757
$ws.Cells.Item(24, 14).Value = 15.76321168151093
$ws.Cells.Item(24, 15).Value = 15.56194183845374
$ws.Cells.Item(25, 2).Value = 10.45047491780055
$ws.Cells.Item(25, 3).Value = 8.694688908328416
$ws.Cells.Item(25, 4).Value = 3.642287087360222
$ws.Cells.Item(25, 6).Value = 17.65435438221946
$ws.Cells.Item(25, 7).Value = 17.71908962591921
$ws.Cells.Item(25, 8).Value = 11.19102974785283
$ws.Cells.Item(25, 9).Value = 15.58279502136999
$ws.Cells.Item(25, 13).Value = 20.85571335814413
$ws.Cells.Item(25, 14).Value = 15.85415585079175
$ws.Cells.Item(25, 15).Value = 15.56941827339475
